$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$setMap = @{
    "N12" = -590
    "J12" = 250
    "L12" = 250
    "H12" = 289.2
    "K17" = 0
    "J17" = 2123
    "H17" = 2123
    "I17" = 0
    "L17" = 6369
    "N17" = -6705
    "I28" = 399.09525
    "M28" = 85.90474999999998
    "H28" = 437.77274
    "K28" = 399.09525
    "I32" = 2168.0476
    "K32" = 2168.0476
    "J32" = 4472
    "L32" = 4472
    "N32" = -5124
    "H32" = 3262.425
    "M32" = -1842.0476
    "J40" = 1500
    "H40" = 1560.9445
    "K40" = 1568.5625
    "N40" = -1850
    "M40" = -1393.5625
    "L40" = 1500
    "I40" = 1568.5625
    "J43" = 0
    "M43" = -224018.89
    "L43" = 0
    "I43" = 224087.89
    "K43" = 224087.89
    "H43" = 224087.89
    "H51" = 6482.1113
    "J51" = 7071.5
    "L51" = 7071.5
    "N51" = -8039.5
    "M53" = -10940
    "I53" = 11577
    "H53" = 6675
    "N53" = -3047
    "L53" = 1773
    "K53" = 11577
    "J53" = 1773
    "J58" = 7499
    "H58" = 2391.2856
    "I58" = 348.2
    "L58" = 22497
    "K58" = 1044.6
    "M58" = -894.5999999999999
    "N58" = -22797
    "N62" = -8247
    "L62" = 6999
    "J62" = 6999
    "M62" = -4197.6665
    "I62" = 4821.6665
    "K62" = 4821.6665
    "H62" = 5217.5454
    "K65" = 24108.3325
    "I65" = 4821.6665
    "N65" = -41235
    "L65" = 34995
    "J65" = 6999
    "M65" = -20988.3325
    "H65" = 5217.5454
    "N86" = -3738.6666
    "H86" = 1382.4286
    "K86" = 1299.75
    "I86" = 1299.75
    "J86" = 1492.6666
    "L86" = 1492.6666
    "M86" = -176.75
    "I88" = 2290.875
    "N88" = -33130.285
    "J88" = 32318.285
    "M88" = -1884.875
    "L88" = 32318.285
    "H88" = 21399.227
    "K88" = 2290.875
    "M89" = -882.75
    "N89" = -18695.333
    "K89" = 6498.75
    "L89" = 7463.333000000001
    "J89" = 1492.6666
    "H89" = 1382.4286
    "I89" = 1299.75
    "K91" = 2290.875
    "L91" = 32318.285
    "M91" = -886.875
    "I91" = 2290.875
    "N91" = -35126.285
    "H91" = 21399.227
    "J91" = 32318.285
    "H92" = 78835.55499999999
    "I92" = 36336.57
    "M92" = -35088.57
    "K92" = 36336.57
    "M111" = 950.2857999999997
    "J111" = 1594
    "H111" = 1075.75
    "L111" = 4782
    "N111" = -10916
    "I111" = 705.5714
    "K111" = 2116.7142
    "J116" = 8134
    "L116" = 8134
    "H116" = 7813.1577
    "K116" = 7372
    "M116" = -3930
    "I116" = 7372
    "N116" = -15018
    "M132" = -1940
    "L132" = 10473.8568
    "N132" = -15533.8568
    "J132" = 3491.2856
    "H132" = 1754.3208
    "K132" = 4470
    "I132" = 1490
    "N135" = -92542.287
    "H135" = 4173.6875
    "K135" = 23588.64
    "M135" = -21053.64
    "I135" = 2620.96
    "L135" = 87472.287
    "J135" = 9719.143
    "I137" = 1641.8
    "H137" = 2302.2954
    "K137" = 4925.4
    "M137" = -2375.4
    "N138" = -23501.75
    "K138" = 4061.5002
    "J138" = 4407.25
    "L138" = 13221.75
    "M138" = 1078.4998
    "I138" = 1353.8334
    "H138" = 3520.7742
}
foreach ($ref in $setMap.Keys) {
    $ws.Range($ref).Value = $setMap[$ref]
}
$delList = @("M17", "N43")
foreach ($ref in $delList) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$setMap = @{
    "L5" = 318.66666
    "I5" = 320.22223
    "H5" = 319.6
    "N5" = -542.66666
    "J5" = 318.66666
    "M5" = -208.22223
    "K5" = 320.22223
    "I32" = 23588.5
    "K32" = 23588.5
    "H32" = 20372.564
    "M32" = -23301.5
    "H45" = 5856.92
    "K45" = 6211.476
    "M45" = -5834.476
    "I45" = 6211.476
    "I61" = 3634.077
    "K61" = 3634.077
    "L61" = 54285.844
    "H61" = 33708.562
    "M61" = -3422.077
    "J61" = 54285.844
    "N61" = -54709.844
    "M74" = -48392.895
    "H74" = 77368
    "K74" = 49266.895
    "I74" = 49266.895
    "H77" = 77368
    "K77" = 246334.475
    "I77" = 49266.895
    "M77" = -241966.475
    "N96" = -61662
    "L96" = 56170
    "J96" = 56170
    "H96" = 56170
    "H122" = 23494.244
    "I122" = 1415.2812
    "K122" = 4245.8436
    "M122" = -1795.8436
    "L132" = 0
    "J132" = 0
    "H132" = 0
    "K132" = 0
    "I132" = 0
    "J136" = 54285.844
    "L136" = 162857.532
    "H136" = 33708.562
    "N136" = -167957.532
    "M136" = -8352.231
    "K136" = 10902.231
    "I136" = 3634.077
}
foreach ($ref in $setMap.Keys) {
    $ws.Range($ref).Value = $setMap[$ref]
}
$delList = @("M132", "N132")
foreach ($ref in $delList) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$setMap = @{
    "N4" = -548.66666
    "K4" = 320.22223
    "J4" = 318.66666
    "I4" = 320.22223
    "H4" = 319.6
    "M4" = -205.22223
    "L4" = 318.66666
    "I7" = 1501.5
    "L7" = 14000
    "J7" = 14000
    "N7" = -14226
    "H7" = 5667.6665
    "M7" = -1388.5
    "K7" = 1501.5
    "M8" = -2425.6667
    "H8" = 2565.6667
    "L8" = 0
    "J8" = 0
    "I8" = 2565.6667
    "K8" = 2565.6667
    "L20" = 2899.5833
    "K20" = 3255.56
    "N20" = -3393.5833
    "J20" = 2899.5833
    "H20" = 3140.1082
    "M20" = -3008.56
    "I20" = 3255.56
    "H22" = 299.25
    "M22" = -177.66666
    "K22" = 350.66666
    "J22" = 145
    "N22" = -491
    "L22" = 145
    "I22" = 350.66666
    "L50" = 0
    "J50" = 0
    "H50" = 0
    "H92" = 17133.666
    "J92" = 17133.666
    "N92" = -22125.666
    "L92" = 17133.666
    "K105" = 1621.3478
    "H105" = 1822.931
    "M105" = 125.6522
    "I105" = 1621.3478
    "L134" = 14903.1432
    "J134" = 4967.7144
    "I134" = 43473.816
    "M134" = -127886.448
    "H134" = 34179.242
    "K134" = 130421.448
    "N134" = -19973.1432
}
foreach ($ref in $setMap.Keys) {
    $ws.Range($ref).Value = $setMap[$ref]
}
$delList = @("N8", "N50")
foreach ($ref in $delList) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$setMap = @{
    "I7" = 53.473682
    "L7" = 98.125
    "J7" = 98.125
    "N7" = -324.125
    "H7" = 66.703705
    "M7" = 59.526318
    "K7" = 53.473682
    "K15" = 0
    "J15" = 17949.5
    "I15" = 0
    "L15" = 17949.5
    "H15" = 17949.5
    "N15" = -18289.5
    "H31" = 2154.762
    "L31" = 2829.3333
    "I31" = 1817.4762
    "M31" = -1522.4762
    "K31" = 1817.4762
    "N31" = -3419.3333
    "J31" = 2829.3333
    "K34" = 1817.4762
    "H34" = 2154.762
    "L34" = 2829.3333
    "I34" = 1817.4762
    "N34" = -3233.3333
    "M34" = -1615.4762
    "J34" = 2829.3333
    "J41" = 48878
    "I41" = 24500
    "M41" = -24072
    "K41" = 24500
    "L41" = 48878
    "H41" = 45395.43
    "N41" = -49734
    "H51" = 21852.889
    "I51" = 16772
    "M51" = -16036
    "K51" = 16772
    "J58" = 27900.8
    "H58" = 20385.715
    "I58" = 1598
    "L58" = 27900.8
    "K58" = 1598
    "M58" = -1395
    "N58" = -28306.8
    "I61" = 16772
    "K61" = 16772
    "H61" = 21852.889
    "M61" = -16424
    "H87" = 64283
    "N87" = -66655
    "I87" = 0
    "J87" = 64283
    "K87" = 0
    "L87" = 64283
    "N90" = -204705
    "H90" = 64283
    "I90" = 0
    "J90" = 64283
    "K90" = 0
    "L90" = 192849
    "K105" = 479.125
    "L105" = 615
    "J105" = 615
    "H105" = 506.3
    "M105" = 1267.875
    "I105" = 479.125
    "N105" = -4109
    "H122" = 1737.4445
    "J122" = 3171.2
    "I122" = 1186
    "L122" = 9513.599999999999
    "K122" = 3558
    "N122" = -14413.6
    "M122" = -1108
    "M132" = -9062.2724
    "H132" = 7400.96
    "K132" = 11592.2724
    "I132" = 3864.0908
    "I134" = 4505.1177
    "M134" = -10980.3531
    "H134" = 4518.1
    "K134" = 13515.3531
    "J136" = 27900.8
    "L136" = 83702.39999999999
    "H136" = 20385.715
    "N136" = -88802.39999999999
    "M136" = -2244
    "K136" = 4794
    "I136" = 1598
}
foreach ($ref in $setMap.Keys) {
    $ws.Range($ref).Value = $setMap[$ref]
}
$delList = @("M15", "M87", "M90")
foreach ($ref in $delList) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$setMap = @{
    "N4" = -46344.5
    "K4" = 1556583.6
    "J4" = 15373.5
    "I4" = 518861.2
    "H4" = 453895.03
    "M4" = -1556471.6
    "L4" = 46120.5
    "L5" = 3000
    "I5" = 445.46667
    "H5" = 480.125
    "N5" = -3224
    "J5" = 1000
    "M5" = -1224.40001
    "K5" = 1336.40001
    "I12" = 800.5
    "M12" = -2228.5
    "K12" = 2401.5
    "N12" = -6309.25
    "J12" = 1987.75
    "L12" = 5963.25
    "H12" = 1818.1428
    "K14" = 5540.5002
    "M14" = -5367.5002
    "H14" = 1846.8334
    "I14" = 1846.8334
    "J38" = 122.28571
    "L38" = 366.85713
    "I38" = 57.833332
    "H38" = 92.53846
    "K38" = 173.499996
    "N38" = -1060.85713
    "M38" = 173.500004
    "M114" = -1295326
    "J114" = 500
    "L114" = 1500
    "N114" = -8008
    "K114" = 1298580
    "I114" = 432860
    "H114" = 324770
    "N129" = -38415.39999999999
    "L129" = 28415.4
    "H129" = 8812.727999999999
    "J129" = 9471.799999999999
    "N135" = -14070
    "H135" = 480.125
    "K135" = 4009.20003
    "M135" = -1474.20003
    "I135" = 445.46667
    "L135" = 9000
    "J135" = 1000
}
foreach ($ref in $setMap.Keys) {
    $ws.Range($ref).Value = $setMap[$ref]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$setMap = @{
    "H22" = 2167
    "M22" = -971.4000000000001
    "K22" = 1500.4
    "J22" = 5500
    "N22" = -6558
    "L22" = 5500
    "I22" = 1500.4
    "L74" = 0
    "J74" = 0
    "H74" = 0
    "L77" = 0
    "H77" = 0
    "J77" = 0
    "J80" = 14491.167
    "L80" = 14491.167
    "H80" = 11087.706
    "I80" = 2919.4
    "N80" = -16487.167
    "K80" = 2919.4
    "M80" = -1921.4
    "J83" = 14491.167
    "K83" = 14597
    "L83" = 72455.83499999999
    "I83" = 2919.4
    "H83" = 11087.706
    "N83" = -82439.83499999999
    "M83" = -9605
    "I102" = 46708.26
    "H102" = 39821.66
    "M102" = -45086.26
    "K102" = 46708.26
    "K113" = 5225.75
    "I113" = 5225.75
    "M113" = -3055.75
    "H113" = 5225.75
    "H122" = 3440.5356
    "J122" = 3875
    "I122" = 2939.2307
    "L122" = 11625
    "K122" = 8817.6921
    "N122" = -16525
    "M122" = -6367.6921
    "M132" = -5228335.4
    "L132" = 18489
    "N132" = -23549
    "J132" = 6163
    "H132" = 1543145.8
    "K132" = 5230865.4
    "I132" = 1743621.8
}
foreach ($ref in $setMap.Keys) {
    $ws.Range($ref).Value = $setMap[$ref]
}
$delList = @("N74", "N77")
foreach ($ref in $delList) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$setMap = @{
    "N4" = -16726
    "K4" = 16672.666
    "J4" = 16500
    "I4" = 16672.666
    "H4" = 16603.6
    "M4" = -16559.666
    "L4" = 16500
    "I7" = 72771.06
    "L7" = 3869.6667
    "J7" = 3869.6667
    "N7" = -4093.6667
    "H7" = 62435.85
    "M7" = -72659.06
    "K7" = 72771.06
    "M16" = -1358.2106
    "I16" = 1528.2106
    "H16" = 2784.1738
    "K16" = 1528.2106
    "N18" = -42412
    "L18" = 42068
    "H18" = 42068.5
    "J18" = 42068
    "J28" = 16500
    "L28" = 16500
    "I28" = 16672.666
    "N28" = -16964
    "M28" = -16440.666
    "H28" = 16603.6
    "K28" = 16672.666
    "N37" = -16714
    "H37" = 16603.6
    "L37" = 16500
    "J37" = 16500
    "I37" = 16672.666
    "M37" = -16565.666
    "K37" = 16672.666
    "K42" = 39999
    "H42" = 39999
    "M42" = -39436
    "I42" = 39999
    "J46" = 7398.6665
    "K46" = 788.6667
    "L46" = 7398.6665
    "M46" = -600.6667
    "I46" = 788.6667
    "N46" = -7774.6665
    "H46" = 4093.6667
    "K49" = 39999
    "I49" = 39999
    "M49" = -39852
    "H49" = 39999
    "I55" = 377.22223
    "M55" = -204.22223
    "N55" = -4970.857
    "H55" = 2235.5625
    "L55" = 4624.857
    "K55" = 377.22223
    "J55" = 4624.857
    "H58" = 7999.6
    "I58" = 3749.75
    "K58" = 3749.75
    "M58" = -3489.75
    "I61" = 3214.7646
    "K61" = 3214.7646
    "H61" = 3214.7646
    "M61" = -3012.7646
    "L100" = 5648
    "H100" = 4195.5386
    "K100" = 2501
    "N100" = -6730
    "J100" = 5648
    "M100" = -1960
    "I100" = 2501
    "K113" = 3214.7646
    "I113" = 3214.7646
    "M113" = -1044.7646
    "H113" = 3214.7646
    "H122" = 4701.4
    "I122" = 4701.4
    "K122" = 14104.2
    "M122" = -11654.2
    "K126" = 218313.18
    "M126" = -215843.18
    "N126" = -16549.0001
    "L126" = 11609.0001
    "J126" = 3869.6667
    "I126" = 72771.06
    "H126" = 62435.85
    "M132" = -6033.875
    "H132" = 3208.652
    "K132" = 8563.875
    "I132" = 2854.625
    "J136" = 6339
    "L136" = 19017
    "H136" = 3208.4092
    "N136" = -24117
    "M136" = -5298.4059
    "K136" = 7848.4059
    "I136" = 2616.1353
}
foreach ($ref in $setMap.Keys) {
    $ws.Range($ref).Value = $setMap[$ref]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$setMap = @{
    "L33" = 37999
    "M33" = -28582.334
    "N33" = -38499
    "I33" = 28832.334
    "H33" = 31124
    "K33" = 28832.334
    "J33" = 37999
    "M36" = -28582.334
    "N36" = -38499
    "I36" = 28832.334
    "J36" = 37999
    "H36" = 31124
    "L36" = 37999
    "K36" = 28832.334
    "I47" = 0
    "H47" = 30069
    "K47" = 0
    "L81" = 7843.5
    "K81" = 2076.1538
    "N81" = -9965.5
    "H81" = 1716.5883
    "I81" = 1038.0769
    "M81" = -1015.1538
    "J81" = 3921.75
    "N84" = -49825.5
    "M84" = -5076.769
    "J84" = 3921.75
    "K84" = 10380.769
    "I84" = 1038.0769
    "L84" = 39217.5
    "H84" = 1716.5883
    "H107" = 487.53333
    "M107" = 923.30769
    "I107" = 332.23077
    "K107" = 996.69231
    "J113" = 2085
    "K113" = 2370.7242
    "I113" = 790.2414
    "N113" = -10595
    "M113" = -200.7242000000001
    "L113" = 6255
    "H113" = 1012.2
    "H122" = 2299.8293
    "J122" = 1750.4286
    "I122" = 2412.9412
    "L122" = 5251.2858
    "K122" = 7238.823600000001
    "N122" = -10151.2858
    "M122" = -4788.823600000001
    "M132" = -86469.00199999999
    "L132" = 89247
    "N132" = -94307
    "J132" = 29749
    "H132" = 29699.4
    "K132" = 88999.00199999999
    "I132" = 29666.334
    "N135" = -86640
    "H135" = 76500
    "L135" = 76500
    "J135" = 76500
    "J136" = 4066
    "L136" = 12198
    "H136" = 3745.12
    "N136" = -17298
    "M136" = -8554.0905
    "K136" = 11104.0905
    "I136" = 3701.3635
}
foreach ($ref in $setMap.Keys) {
    $ws.Range($ref).Value = $setMap[$ref]
}
$delList = @("M47")
foreach ($ref in $delList) {
    $ws.Range($ref).ClearContents()
}
